# Update the division-problem answers in the single table of the document.
# The table has 5 "data" rows (1, 5, 9, 13, 17) each holding 5 problems.
# We rewrite each cell's text in place with Cell.Range.Text so the
# existing run formatting (rFonts/sz) defined on the paragraph's run is
# preserved by the runtime, matching how Word normally replaces cell text.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "77÷4=19, 1"
$t.Cell(1,2).Range.Text = "95÷3=31, 2"
$t.Cell(1,3).Range.Text = "81÷8=10, 1"
$t.Cell(1,4).Range.Text = "77÷8=9, 5"
$t.Cell(1,5).Range.Text = "94÷5=18, 4"

$t.Cell(5,1).Range.Text = "30÷7=4, 2"
$t.Cell(5,2).Range.Text = "32÷3=10, 2"
$t.Cell(5,3).Range.Text = "31÷3=10, 1"
$t.Cell(5,4).Range.Text = "20÷9=2, 2"
$t.Cell(5,5).Range.Text = "49÷4=12, 1"

$t.Cell(9,1).Range.Text = "24÷7=3, 3"
$t.Cell(9,2).Range.Text = "38÷9=4, 2"
$t.Cell(9,3).Range.Text = "19÷7=2, 5"
$t.Cell(9,4).Range.Text = "96÷2=48, 0"
$t.Cell(9,5).Range.Text = "50÷9=5, 5"

$t.Cell(13,1).Range.Text = "26÷6=4, 2"
$t.Cell(13,2).Range.Text = "45÷6=7, 3"
$t.Cell(13,3).Range.Text = "41÷7=5, 6"
$t.Cell(13,4).Range.Text = "25÷2=12, 1"
$t.Cell(13,5).Range.Text = "17÷6=2, 5"

$t.Cell(17,1).Range.Text = "94÷5=18, 4"
$t.Cell(17,2).Range.Text = "19÷5=3, 4"
$t.Cell(17,3).Range.Text = "28÷9=3, 1"
$t.Cell(17,4).Range.Text = "19÷7=2, 5"
$t.Cell(17,5).Range.Text = "19÷2=9, 1"
